# Cthrc1-Fzd5.xlsx was regenerated with a new TPM run.
# The ligand-receptor table now reports all 3x3 sender/target cluster
# combinations (ECs, FAPs, MuSCs) instead of only 2 of the 3 target
# clusters per sender, growing the data block from 6 rows to 9 rows
# (sheet dimension A1:T7 -> A1:T10). Column A/D values, and every numeric
# statistic in columns E:T, are updated to match the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One entry per output data row (sheet rows 2-10), in final top-to-bottom
# order: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric statistic columns E..T.
$rows = @(
    @{ A="ECs";   D="ECs";   E=1; F=0.3333333333333333; G=0.01838633333333333; H=0.055159;   I=0.003339500866342531;  J=0.003339500866342531;  K=3; L=1; M=2.133443333333334; N=6.40033;   O=0.2605947899689859; P=0.2605947899689859; Q=0.03922620027444445; R=0.35303580247;      S=0.0008702565268657785; T=0.0008702565268657784 },
    @{ A="ECs";   D="FAPs";  E=1; F=0.3333333333333333; G=0.01838633333333333; H=0.055159;   I=0.003339500866342531;  J=0.003339500866342531;  K=3; L=1; M=4.264793333333333; N=12.79438;  O=0.5209338844846115; P=0.5209338844846116; Q=0.07841391182444445; R=0.70572520642;      S=0.00173965915854354;   T=0.00173965915854354   },
    @{ A="ECs";   D="MuSCs"; E=1; F=0.3333333333333333; G=0.01838633333333333; H=0.055159;   I=0.003339500866342531;  J=0.003339500866342531;  K=3; L=1; M=1.788586;          N=5.365758;  O=0.2184713255464024; P=0.2184713255464024; Q=0.03288553839133333; R=0.295969845522;     S=0.0007295851809332119; T=0.000729585180933212  },
    @{ A="FAPs";  D="ECs";   E=3; F=1;                  G=5.449420333333333;   H=16.348261;  I=0.9897755900704113;    J=0.9897755900704112;    K=3; L=1; M=2.133443333333334; N=6.40033;   O=0.2605947899689859; P=0.2605947899689859; Q=11.62602948068111;   R=104.63426532613;    S=0.257930362010828;     T=0.2579303620108279    },
    @{ A="FAPs";  D="FAPs";  E=3; F=1;                  G=5.449420333333333;   H=16.348261;  I=0.9897755900704113;    J=0.9897755900704112;    K=3; L=1; M=4.264793333333333; N=12.79438;  O=0.5209338844846115; P=0.5209338844846116; Q=23.24065150813111;   R=209.16586357318;    S=0.5156076429034279;    T=0.5156076429034279    },
    @{ A="FAPs";  D="MuSCs"; E=3; F=1;                  G=5.449420333333333;   H=16.348261;  I=0.9897755900704113;    J=0.9897755900704112;    K=3; L=1; M=1.788586;          N=5.365758;  O=0.2184713255464024; P=0.2184713255464024; Q=9.746756916315332;   R=87.72081224683799;  S=0.2162375851561553;    T=0.2162375851561554    },
    @{ A="MuSCs"; D="ECs";   E=1; F=0.3333333333333333; G=0.03790633333333333; H=0.113719;   I=0.006884909063246366;  J=0.006884909063246365;  K=3; L=1; M=2.133443333333334; N=6.40033;   O=0.2605947899689859; P=0.2605947899689859; Q=0.08087101414111113; R=0.72783912727;      S=0.001794171431292255;  T=0.001794171431292254  },
    @{ A="MuSCs"; D="FAPs";  E=1; F=0.3333333333333333; G=0.03790633333333333; H=0.113719;   I=0.006884909063246366;  J=0.006884909063246365;  K=3; L=1; M=4.264793333333333; N=12.79438;  O=0.5209338844846115; P=0.5209338844846116; Q=0.1616626776911111;  R=1.45496409922;      S=0.003586582422640238;  T=0.003586582422640238  },
    @{ A="MuSCs"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.03790633333333333; H=0.113719;   I=0.006884909063246366;  J=0.006884909063246365;  K=3; L=1; M=1.788586;          N=5.365758;  O=0.2184713255464024; P=0.2184713255464024; Q=0.06779873711133333; R=0.610188634002;     S=0.001504155209313873;  T=0.001504155209313873  }
)

$ligand = "Cthrc1"
$receptor = "Fzd5"
$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rows) {
    $row["B"] = $ligand
    $row["C"] = $receptor

    for ($c = 0; $c -lt $colOrder.Count; $c++) {
        $col = $colOrder[$c]
        $ws.Cells.Item($r, $c + 1).Value = $row[$col]
    }

    $r = $r + 1
}
